$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.899.72'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').Value = '3.263.08'
$ws.Range('E3').Value = '  +2.64%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '396.57'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.15'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.586'
$ws.Range('E7').Value = '  +6.40%  '
$ws.Range('D8').Value = '3.261.47'
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.628'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '39.35'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('E12').Value = '  +10.52%  '
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').Value = '3.761.63'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.44'
$ws.Range('E15').Value = '  +4.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.20'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '3.248.62'
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('E18').Value = '  -2.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.77'
$ws.Range('E19').Value = '  +2.10%  '
$ws.Range('D20').Value = '56.802.79'
$ws.Range('E20').Value = '  +4.03%  '
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('E22').Value = '  +8.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.01'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '295.98'
$ws.Range('E24').Value = '  +7.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.47'
$ws.Range('E25').Value = '  +2.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.18'
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.21'
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.68'
$ws.Range('E29').Value = '  -3.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.30'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.169'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.28'
$ws.Range('E33').Value = '  +2.10%  '
$ws.Range('E34').Value = '  -3.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.00'
$ws.Range('E35').Value = '  +8.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0490'
$ws.Range('E36').Value = '  -3.70%  '
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.47'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.49'
$ws.Range('E40').Value = '  -3.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.94'
$ws.Range('E41').Value = '  +2.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '137.69'
$ws.Range('E42').Value = '  +5.05%  '
$ws.Range('E43').Value = '  +3.82%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.02'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.91'
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.15'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.25'
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('E49').Value = '  +3.32%  '
$ws.Range('D50').Value = '2.158.87'
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('E51').Value = '  -4.95%  '
